# Add "Baseline" analyses folder column to the subjects sheet.
#
# Original layout : A=name B=folder(int path) C=group D=int E=grade F=age G=sub H=raw_data
# Target layout   : A=name B=raw_data C=folder_baseline D=folder_int E=group F=int G=grade H=age I=sub
#
# Strategy:
#   1. Move column H (raw_data) to become the new column B (shifts B..G -> C..H).
#   2. Insert a brand-new blank column at C for "folder_baseline" (shifts old C..H -> D..I).
#   3. Rename the old "folder" header to "folder_int".
#   4. Populate the new "folder_baseline" column header + the per-subject paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relocate raw_data (H) so it immediately follows name (A).
$ws.Columns.Item(8).Cut()
$ws.Columns.Item(2).Insert()

# 2. Make room for the new folder_baseline column right after raw_data.
$ws.Columns.Item(3).Insert()

# 3. Header renames / additions.
$ws.Range("D1").Value = "folder_int"
$ws.Range("C1").Value = "folder_baseline"

# 4. Fill in the per-subject Baseline derivative paths (mirrors folder_int, with sub-NN\Baseline).
$ws.Range("C2").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-01\Baseline"
$ws.Range("C3").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-02\Baseline"
$ws.Range("C4").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-03\Baseline"
$ws.Range("C5").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-04\Baseline"
$ws.Range("C6").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-05\Baseline"
$ws.Range("C7").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-06\Baseline"
$ws.Range("C8").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-07\Baseline"
$ws.Range("C9").Value  = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-08\Baseline"
$ws.Range("C10").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-09\Baseline"
$ws.Range("C11").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-10\Baseline"
$ws.Range("C12").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-11\Baseline"
$ws.Range("C13").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-12\Baseline"
$ws.Range("C14").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-13\Baseline"
$ws.Range("C15").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-14\Baseline"
$ws.Range("C16").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-15\Baseline"
$ws.Range("C17").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-16\Baseline"
$ws.Range("C18").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-17\Baseline"
$ws.Range("C19").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-18\Baseline"
$ws.Range("C20").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-19\Baseline"
$ws.Range("C21").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-20\Baseline"
$ws.Range("C22").Value = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-21\Baseline"

# Update the view: clear the scrolled top-left cell and move the active selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C23").Select()
